# Auto-generated Excel COM-interop script
# Applies per-cell value updates (and a few cell deletions) to the
# leve-profit calculator sheets, matching the upstream market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4999.5
$ws.Range("H72").Value = 4999.5
$ws.Range("H86").Value = 5400.5
$ws.Range("I86").Value = 4471.7144
$ws.Range("K86").Value = 4471.7144
$ws.Range("M86").Value = -3348.7144
$ws.Range("H88").Value = 2138.3809
$ws.Range("J88").Value = 2356
$ws.Range("L88").Value = 2356
$ws.Range("N88").Value = -3168
$ws.Range("H89").Value = 5400.5
$ws.Range("I89").Value = 4471.7144
$ws.Range("K89").Value = 22358.572
$ws.Range("M89").Value = -16742.572
$ws.Range("H91").Value = 2138.3809
$ws.Range("J91").Value = 2356
$ws.Range("L91").Value = 2356
$ws.Range("N91").Value = -5164
$ws.Range("H127").Value = 850.5
$ws.Range("I127").Value = 800.6667
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 2402.0001
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = 2557.9999
$ws.Range("N127").Value = -12920
$ws.Range("H141").Value = 21750
$ws.Range("I141").Value = 21750
$ws.Range("K141").Value = 65250
$ws.Range("M141").Value = -60070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 5018.75
$ws.Range("I3").Value = 1037.5
$ws.Range("K3").Value = 1037.5
$ws.Range("M3").Value = -922.5
$ws.Range("H24").Value = 77663.75
$ws.Range("J24").Value = 77663.75
$ws.Range("L24").Value = 77663.75
$ws.Range("N24").Value = -78411.75
$ws.Range("H88").Value = 3137.4
$ws.Range("I88").Value = 1770.25
$ws.Range("J88").Value = 4699.857
$ws.Range("K88").Value = 1770.25
$ws.Range("L88").Value = 4699.857
$ws.Range("M88").Value = -1364.25
$ws.Range("N88").Value = -5511.857
$ws.Range("H91").Value = 3137.4
$ws.Range("I91").Value = 1770.25
$ws.Range("J91").Value = 4699.857
$ws.Range("K91").Value = 1770.25
$ws.Range("L91").Value = 4699.857
$ws.Range("M91").Value = -366.25
$ws.Range("N91").Value = -7507.857
$ws.Range("H92").Value = 46750.25
$ws.Range("J92").Value = 46750.25
$ws.Range("L92").Value = 46750.25
$ws.Range("N92").Value = -51742.25
$ws.Range("H96").Value = 22999
$ws.Range("J96").Value = 22999
$ws.Range("L96").Value = 22999
$ws.Range("N96").Value = -28491
$ws.Range("H100").Value = 77663.75
$ws.Range("J100").Value = 77663.75
$ws.Range("L100").Value = 77663.75
$ws.Range("N100").Value = -79827.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7499.25
$ws.Range("I86").Value = 7499.25
$ws.Range("K86").Value = 7499.25
$ws.Range("M86").Value = -6376.25
$ws.Range("H89").Value = 7499.25
$ws.Range("I89").Value = 7499.25
$ws.Range("K89").Value = 37496.25
$ws.Range("M89").Value = -31880.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1256.6364
$ws.Range("I122").Value = 1855
$ws.Range("K122").Value = 5565
$ws.Range("M122").Value = -3115
$ws.Range("H132").Value = 9717.6
$ws.Range("I132").Value = 8423.5
$ws.Range("J132").Value = 14894
$ws.Range("K132").Value = 25270.5
$ws.Range("L132").Value = 44682
$ws.Range("M132").Value = -22740.5
$ws.Range("N132").Value = -49742

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 543
$ws.Range("I7").Value = 543
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1629
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1517
$ws.Range("N7").ClearContents()
$ws.Range("H92").Value = 744.8125
$ws.Range("I92").Value = 520.6667
$ws.Range("J92").Value = 879.3
$ws.Range("K92").Value = 1562.0001
$ws.Range("L92").Value = 2637.9
$ws.Range("M92").Value = -314.0001
$ws.Range("N92").Value = -5133.9
$ws.Range("H107").Value = 714.2
$ws.Range("I107").Value = 599.8
$ws.Range("J107").Value = 828.6
$ws.Range("K107").Value = 1799.4
$ws.Range("L107").Value = 2485.8
$ws.Range("M107").Value = 120.6000000000001
$ws.Range("N107").Value = -6325.8
$ws.Range("H134").Value = 250003000
$ws.Range("I134").Value = 250003000
$ws.Range("K134").Value = 750009000
$ws.Range("M134").Value = -750003930
$ws.Range("H137").Value = 1725.5714
$ws.Range("I137").Value = 733.6667
$ws.Range("J137").Value = 2469.5
$ws.Range("K137").Value = 2201.0001
$ws.Range("L137").Value = 7408.5
$ws.Range("M137").Value = 2898.9999
$ws.Range("N137").Value = -17608.5
$ws.Range("H139").Value = 1681.5
$ws.Range("I139").Value = 1681.5
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5044.5
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 95.5
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1499.4
$ws.Range("I80").Value = 1649.75
$ws.Range("J80").Value = 898
$ws.Range("K80").Value = 1649.75
$ws.Range("L80").Value = 898
$ws.Range("M80").Value = -651.75
$ws.Range("N80").Value = -2894
$ws.Range("H83").Value = 1499.4
$ws.Range("I83").Value = 1649.75
$ws.Range("J83").Value = 898
$ws.Range("K83").Value = 8248.75
$ws.Range("L83").Value = 4490
$ws.Range("M83").Value = -3256.75
$ws.Range("N83").Value = -14474
$ws.Range("H97").Value = 1036
$ws.Range("I97").Value = 769.53845
$ws.Range("K97").Value = 769.53845
$ws.Range("M97").Value = -273.53845
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 1763
$ws.Range("J122").Value = 1758
$ws.Range("L122").Value = 5274
$ws.Range("N122").Value = -10174

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2478.2
$ws.Range("I40").Value = 2478.2
$ws.Range("K40").Value = 2478.2
$ws.Range("M40").Value = -2342.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1759.6923
$ws.Range("I96").Value = 1558.7
$ws.Range("J96").Value = 2429.6667
$ws.Range("K96").Value = 1558.7
$ws.Range("L96").Value = 2429.6667
$ws.Range("M96").Value = -185.7
$ws.Range("N96").Value = -5175.6667
$ws.Range("H113").Value = 548.5
$ws.Range("I113").Value = 548.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1645.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 524.5
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 1681
$ws.Range("I136").Value = 1167
$ws.Range("J136").Value = 4336.6665
$ws.Range("K136").Value = 3501
$ws.Range("L136").Value = 13009.9995
$ws.Range("M136").Value = -951
$ws.Range("N136").Value = -18109.9995
